$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 (shifts rows 8-14 down to 9-15)
$ws.Rows.Item(8).Insert()

# Fill in the new row's data
$ws.Cells.Item(8, 1).Value = "mistral_7b_instruct_v2"
$ws.Cells.Item(8, 2).Value = 91.22935779816514
$ws.Cells.Item(8, 3).Value = 38.86238532110092
$ws.Cells.Item(8, 4).Value = 47.52293577981651
$ws.Cells.Item(8, 5).Value = 39.26605504587156
$ws.Cells.Item(8, 6).Value = 5.688073394495413
$ws.Cells.Item(8, 7).Value = 3.926605504587156
$ws.Cells.Item(8, 8).Value = 2.422018348623853
$ws.Cells.Item(8, 9).Value = 2.642201834862385
